$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-53 with new date/value pairs
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 0.7960900427496256
$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2.427463180990628
$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = -0.2156301208428602
$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = -9.679894544238294
$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = -0.8195522748085864
$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2.274906658561562
$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 7.912074126222677
$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 1.840273951972307
$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2.530590017523977
$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 1.726328558848905
$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 1.446992732417414
$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = -1.803419669653366
$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2.429152676597084
$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 0.1648657400322122
$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 1.153911995371715
$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 0.9978838984052771
$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 1.7958205029814
$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = -0.6421078521012049
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = -0.4018974900578201
$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 1.682996656587392
$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 1.044515111522799
$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 1.750939777838624
$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = -1.040686374073175
$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = -0.2631342424776193
$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 0.8
$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = -0.8578072826671388
$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 1
$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 1.83469776744667
$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = -1.25257285039946
$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 0.9704313550711277
$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = -3.1
$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = -20.33630385521235
$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 18.14253003720006
$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 4.544557841334466
$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 1.849039957006937
$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 0.5397909189180297
$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = -0.9815688985294599
$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 4.784172743600905
$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = -2.12351799362105
$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 0.2867860097325945
$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 1.982042055103904
$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = -0.9572753870817365
$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 0.4216256340467766
$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = -1.083093549703378
$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = -0.756372051364167
$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = -1.594471108675521
$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 1.08522880488573
$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = -0.150193889177217
$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = -1.874498177919662
$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = -2.24094337053269
$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 3.234504069499366
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = -0.1262090121580286

# Remove old trailing rows 54-73 (shift cells up, deleting the rows entirely)
$ws.Range("A54:B73").EntireRow.Delete()

Write-Output ("Final UsedRange: " + $ws.UsedRange.Address())
